$wb = $excel.ActiveWorkbook

# Sheet "展览" - rows 5-8 in column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 3246
$wsExhibit.Range("F6").Value = 329
$wsExhibit.Range("F7").Value = 11
$wsExhibit.Range("F8").Value = 415

# Sheet "全部类型" - rows 5,6,9,10 in column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 3246
$wsAll.Range("F6").Value = 329
$wsAll.Range("F9").Value = 11
$wsAll.Range("F10").Value = 415
